$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Supplier" column (K) — header + a supplier id (2) for every data row.
$ws.Range("K1").Value = "Supplier"
$ws.Range("K2:K7").Value = 2

# Row 2 becomes slightly shorter in the edited workbook.
$ws.Rows.Item(2).RowHeight = 15

# Put the selection on the newly added cell, matching the author's edit.
[void]$ws.Range("K7").Select()
